# Session 4: Greedy algorithms
# - Record a "10" test score for the student in the "Session 4 (Greedy)"
#   column (E4), which also ripples into the Mark/J4 formula recalculation.
# - Add the grader's remark "Good!" in the feedback row (E5) under that
#   same session column.
# - Move the active selection from the old Session-4 column (E5:E12) to the
#   next one over (F5:F12), matching where the user was working next.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 10
$ws.Range("E5").Value = "Good!"

$ws.Range("F5:F12").Select()
